$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14's "Descripción" note was an explicit-but-empty value; the refreshed
# export no longer emits it at all, so the cell becomes fully blank.
$ws.Range("D14").ClearContents()

# --- New attendance entries recorded on 2025-03-25 ---
# Dates are leading-quoted so Excel stores them as literal text (matching
# every other row in the log) instead of auto-converting them to a date
# serial number; re-applying the plain "Normal" style afterwards strips the
# transient quote-prefix formatting that assignment leaves behind.

# Row 15: Salida at 00:47:14 (this row has no Descripción cell at all).
$ws.Range("A15").Value = "'2025-03-25"
$ws.Range("A15").Style = "Normal"
$ws.Range("B15").Value = "00:47:14"
$ws.Range("C15").Value = "Salida"

# Row 16: Salida at 00:48:23, with an explicit-but-empty Descripción.
$ws.Range("A16").Value = "'2025-03-25"
$ws.Range("A16").Style = "Normal"
$ws.Range("B16").Value = "00:48:23"
$ws.Range("C16").Value = "Salida"
$ws.Range("D16").Value = "'"
$ws.Range("D16").Style = "Normal"

# Row 17: Salida at 00:49:14, with an explicit-but-empty Descripción.
$ws.Range("A17").Value = "'2025-03-25"
$ws.Range("A17").Style = "Normal"
$ws.Range("B17").Value = "00:49:14"
$ws.Range("C17").Value = "Salida"
$ws.Range("D17").Value = "'"
$ws.Range("D17").Style = "Normal"

# Row 18: Entrada at 00:49:36, with an explicit-but-empty Descripción.
$ws.Range("A18").Value = "'2025-03-25"
$ws.Range("A18").Style = "Normal"
$ws.Range("B18").Value = "00:49:36"
$ws.Range("C18").Value = "Entrada"
$ws.Range("D18").Value = "'"
$ws.Range("D18").Style = "Normal"

# Row 19: Salida at 00:49:42, with an explicit-but-empty Descripción.
$ws.Range("A19").Value = "'2025-03-25"
$ws.Range("A19").Style = "Normal"
$ws.Range("B19").Value = "00:49:42"
$ws.Range("C19").Value = "Salida"
$ws.Range("D19").Value = "'"
$ws.Range("D19").Style = "Normal"
